$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 64999.668
$ws.Range("J3").Value = 64999.668
$ws.Range("L3").Value = 64999.668
$ws.Range("N3").Value = -65227.668
$ws.Range("H4").Value = 274.2
$ws.Range("I4").Value = 274.2
$ws.Range("K4").Value = 274.2
$ws.Range("M4").Value = -160.2
$ws.Range("H18").Value = 62500132
$ws.Range("I18").Value = 62500132
$ws.Range("K18").Value = 62500132
$ws.Range("M18").Value = -62499848
$ws.Range("I43").Value = 1024.75
$ws.Range("J43").Value = 1599.3334
$ws.Range("K43").Value = 1024.75
$ws.Range("L43").Value = 1599.3334
$ws.Range("M43").Value = -955.75
$ws.Range("N43").Value = -1737.3334
$ws.Range("H51").Value = 4222.222
$ws.Range("J51").Value = 4400
$ws.Range("L51").Value = 4400
$ws.Range("N51").Value = -5368
$ws.Range("H53").Value = 794.36365
$ws.Range("I53").Value = 990.8570999999999
$ws.Range("K53").Value = 990.8570999999999
$ws.Range("M53").Value = -353.8570999999999
$ws.Range("H58").Value = 1744.2307
$ws.Range("I58").Value = 243.18182
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 729.5454599999999
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -579.5454599999999
$ws.Range("N58").Value = -30300
$ws.Range("H80").Value = 950.6316
$ws.Range("I80").Value = 692.7646999999999
$ws.Range("J80").Value = 1159.381
$ws.Range("K80").Value = 2078.2941
$ws.Range("L80").Value = 3478.143
$ws.Range("M80").Value = -1080.2941
$ws.Range("N80").Value = -5474.143
$ws.Range("H83").Value = 950.6316
$ws.Range("I83").Value = 692.7646999999999
$ws.Range("J83").Value = 1159.381
$ws.Range("K83").Value = 6234.882299999999
$ws.Range("L83").Value = 10434.429
$ws.Range("M83").Value = -1242.882299999999
$ws.Range("N83").Value = -20418.429
$ws.Range("H102").Value = 64999.668
$ws.Range("J102").Value = 64999.668
$ws.Range("L102").Value = 64999.668
$ws.Range("N102").Value = -71489.66800000001
$ws.Range("H132").Value = 4575.62
$ws.Range("I132").Value = 2628.4443
$ws.Range("J132").Value = 12242.625
$ws.Range("K132").Value = 7885.3329
$ws.Range("L132").Value = 36727.875
$ws.Range("M132").Value = -5355.3329
$ws.Range("N132").Value = -41787.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 4190.7144
$ws.Range("I21").Value = 5102.5713
$ws.Range("J21").Value = 3278.8572
$ws.Range("K21").Value = 5102.5713
$ws.Range("L21").Value = 3278.8572
$ws.Range("M21").Value = -4728.5713
$ws.Range("N21").Value = -4026.8572
$ws.Range("H30").Value = 2947.7273
$ws.Range("I30").Value = 1292
$ws.Range("J30").Value = 3893.8572
$ws.Range("K30").Value = 1292
$ws.Range("L30").Value = 3893.8572
$ws.Range("M30").Value = -1142
$ws.Range("N30").Value = -4193.8572
$ws.Range("H32").Value = 4067.9011
$ws.Range("I32").Value = 3575.8838
$ws.Range("J32").Value = 12530.6
$ws.Range("K32").Value = 3575.8838
$ws.Range("L32").Value = 12530.6
$ws.Range("M32").Value = -3288.8838
$ws.Range("N32").Value = -13104.6
$ws.Range("H45").Value = 10876.056
$ws.Range("I45").Value = 12447.692
$ws.Range("J45").Value = 6789.8
$ws.Range("K45").Value = 12447.692
$ws.Range("L45").Value = 6789.8
$ws.Range("M45").Value = -12070.692
$ws.Range("N45").Value = -7543.8
$ws.Range("H105").Value = 87500
$ws.Range("J105").Value = 87500
$ws.Range("L105").Value = 87500
$ws.Range("N105").Value = -94488
$ws.Range("H122").Value = 18921.05
$ws.Range("I122").Value = 2581.3845
$ws.Range("K122").Value = 7744.1535
$ws.Range("M122").Value = -5294.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1344.0625
$ws.Range("I86").Value = 1399.8096
$ws.Range("J86").Value = 1237.6364
$ws.Range("K86").Value = 1399.8096
$ws.Range("L86").Value = 1237.6364
$ws.Range("M86").Value = -276.8096
$ws.Range("N86").Value = -3483.6364
$ws.Range("H89").Value = 1344.0625
$ws.Range("I89").Value = 1399.8096
$ws.Range("J89").Value = 1237.6364
$ws.Range("K89").Value = 6999.048000000001
$ws.Range("L89").Value = 6188.182000000001
$ws.Range("M89").Value = -1383.048000000001
$ws.Range("N89").Value = -17420.182
$ws.Range("H97").Value = 1990
$ws.Range("I97").Value = 1990
$ws.Range("K97").Value = 1990
$ws.Range("M97").Value = -999
$ws.Range("H99").Value = 1210.5385
$ws.Range("I99").Value = 1277.091
$ws.Range("K99").Value = 1277.091
$ws.Range("M99").Value = 220.9090000000001
$ws.Range("H134").Value = 4420.1377
$ws.Range("I134").Value = 3196.5217
$ws.Range("K134").Value = 9589.5651
$ws.Range("M134").Value = -7054.5651

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1955.3455
$ws.Range("I31").Value = 1955.3077
$ws.Range("K31").Value = 1955.3077
$ws.Range("M31").Value = -1660.3077
$ws.Range("H32").Value = 6712.8
$ws.Range("I32").Value = 5104.8335
$ws.Range("J32").Value = 9124.75
$ws.Range("K32").Value = 5104.8335
$ws.Range("L32").Value = 9124.75
$ws.Range("M32").Value = -4788.8335
$ws.Range("N32").Value = -9756.75
$ws.Range("H34").Value = 1955.3455
$ws.Range("I34").Value = 1955.3077
$ws.Range("K34").Value = 1955.3077
$ws.Range("M34").Value = -1753.3077
$ws.Range("H81").Value = 69998.5
$ws.Range("J81").Value = 69998.5
$ws.Range("L81").Value = 69998.5
$ws.Range("N81").Value = -71994.5
$ws.Range("H84").Value = 69998.5
$ws.Range("J84").Value = 69998.5
$ws.Range("L84").Value = 209995.5
$ws.Range("N84").Value = -219979.5
$ws.Range("H134").Value = 2877.0789
$ws.Range("I134").Value = 3004.1292
$ws.Range("J134").Value = 2314.4285
$ws.Range("K134").Value = 9012.3876
$ws.Range("L134").Value = 6943.2855
$ws.Range("M134").Value = -6477.3876
$ws.Range("N134").Value = -12013.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 274888.53
$ws.Range("I4").Value = 481182.44
$ws.Range("J4").Value = 4127.8125
$ws.Range("K4").Value = 1443547.32
$ws.Range("L4").Value = 12383.4375
$ws.Range("M4").Value = -1443435.32
$ws.Range("N4").Value = -12607.4375
$ws.Range("H33").Value = 305.76923
$ws.Range("J33").Value = 374.1111
$ws.Range("L33").Value = 2244.6666
$ws.Range("N33").Value = -2810.6666
$ws.Range("H97").Value = 662.4545000000001
$ws.Range("J97").Value = 772.25
$ws.Range("L97").Value = 2316.75
$ws.Range("N97").Value = -3308.75
$ws.Range("H131").Value = 3679844.5
$ws.Range("J131").Value = 3856.1428
$ws.Range("L131").Value = 11568.4284
$ws.Range("N131").Value = -21648.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 258243
$ws.Range("J3").Value = 9591.6
$ws.Range("L3").Value = 9591.6
$ws.Range("N3").Value = -9823.6
$ws.Range("H12").Value = 2932.6667
$ws.Range("I12").Value = 2932.6667
$ws.Range("K12").Value = 2932.6667
$ws.Range("M12").Value = -2792.6667
$ws.Range("H19").Value = 6500
$ws.Range("I19").Value = 7000
$ws.Range("K19").Value = 7000
$ws.Range("M19").Value = -6712
$ws.Range("H43").Value = 12800
$ws.Range("I43").Value = 6000
$ws.Range("K43").Value = 6000
$ws.Range("M43").Value = -5849
$ws.Range("H57").Value = 27438.666
$ws.Range("J57").Value = 30870.5
$ws.Range("L57").Value = 30870.5
$ws.Range("N57").Value = -32510.5
$ws.Range("H70").Value = 15371.143
$ws.Range("J70").Value = 10444
$ws.Range("L70").Value = 10444
$ws.Range("N70").Value = -10984
$ws.Range("H73").Value = 15371.143
$ws.Range("J73").Value = 10444
$ws.Range("L73").Value = 10444
$ws.Range("N73").Value = -12316
$ws.Range("H80").Value = 48002720
$ws.Range("I80").Value = 80001016
$ws.Range("J80").Value = 5280.25
$ws.Range("K80").Value = 80001016
$ws.Range("L80").Value = 5280.25
$ws.Range("M80").Value = -80000018
$ws.Range("N80").Value = -7276.25
$ws.Range("H83").Value = 48002720
$ws.Range("I83").Value = 80001016
$ws.Range("J83").Value = 5280.25
$ws.Range("K83").Value = 400005080
$ws.Range("L83").Value = 26401.25
$ws.Range("M83").Value = -400000088
$ws.Range("N83").Value = -36385.25
$ws.Range("H122").Value = 3291.5
$ws.Range("I122").Value = 2801.875
$ws.Range("K122").Value = 8405.625
$ws.Range("M122").Value = -5955.625
$ws.Range("H132").Value = 6305.25
$ws.Range("I132").Value = 3994.2354
$ws.Range("K132").Value = 11982.7062
$ws.Range("M132").Value = -9452.706200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 12248.75
$ws.Range("I26").Value = 12999
$ws.Range("K26").Value = 12999
$ws.Range("M26").Value = -12704
$ws.Range("H40").Value = 6531.2
$ws.Range("I40").Value = 5759.5713
$ws.Range("K40").Value = 5759.5713
$ws.Range("M40").Value = -5623.5713
$ws.Range("H106").Value = 21165.834
$ws.Range("J106").Value = 21165.834
$ws.Range("L106").Value = 21165.834
$ws.Range("N106").Value = -23689.834
$ws.Range("H136").Value = 3832228.2
$ws.Range("I136").Value = 5001714.5
$ws.Range("J136").Value = 4817.909
$ws.Range("K136").Value = 15005143.5
$ws.Range("L136").Value = 14453.727
$ws.Range("M136").Value = -15002593.5
$ws.Range("N136").Value = -19553.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 63399.57
$ws.Range("J26").Value = 56999.25
$ws.Range("L26").Value = 56999.25
$ws.Range("N26").Value = -57585.25
$ws.Range("H104").Value = 27589.445
$ws.Range("J104").Value = 27589.445
$ws.Range("L104").Value = 27589.445
$ws.Range("N104").Value = -34577.445
$ws.Range("H136").Value = 1618.4769
$ws.Range("I136").Value = 1621.1666
$ws.Range("K136").Value = 4863.4998
$ws.Range("M136").Value = -2313.4998
